# Resort the worksheets: move "总计" (the summary/total sheet, currently
# second) in front of "2022-Q2" (the per-fund detail sheet, currently
# first), so the tab order becomes: 总计, 2022-Q2.
$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

$totalSheet.Move($q2Sheet)
